# feat: add 2022-Q3 data
#
# Before: sheets = [ "总计", "2021-Q4" ]
# After:  sheets = [ "总计", "2022-Q3", "2021-Q4" ]
#   - "2021-Q4" (old) is duplicated -> the duplicate keeps the name "2021-Q4"
#     and retains the old Q4 numbers (this becomes the new, 3rd sheet).
#   - the original "2021-Q4" sheet is renamed to "2022-Q3" and its data is
#     overwritten with the new Q3 numbers.
#   - the "总计" (totals) sheet gets a new row for 2022-Q3 (inserted above
#     the existing 2021-Q4 total row) and its own 2021-Q4 row is preserved
#     unchanged below it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a value to be written as TEXT (inlineStr/shared-string) even
# when it looks like a number ("45.98", "001481", ...), without leaving a
# lasting custom number-format on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Step 1: duplicate the "2021-Q4" sheet so the old data survives on its own
# tab, placed immediately after the source sheet.
# ---------------------------------------------------------------------------
$q4Old = $wb.Worksheets.Item("2021-Q4")
$q4Old.Copy($null, $q4Old)

# The copy lands right after $q4Old and is auto-named "2021-Q4 (2)". Rename
# the original sheet (still holding the old Q4 numbers for now) to its new
# identity FIRST, freeing up the "2021-Q4" name for the duplicate.
$q3 = $q4Old
$q3.Name = "2022-Q3"

# Now put the "2021-Q4" name on the duplicate (the one that keeps the
# historical name/data).
$q4Copy = $wb.Worksheets.Item($q3.Index + 1)
$q4Copy.Name = "2021-Q4"

# ---------------------------------------------------------------------------
# Step 2: overwrite the (renamed) sheet's data rows with the 2022-Q3 figures.
# Columns A (index) and B (fund code) are identical to the old data, so they
# are left untouched; C-H get new values.
# ---------------------------------------------------------------------------
Set-TextValue $q3.Cells.Item(2, 3) "华宝油气（QDII）美元"
Set-TextValue $q3.Cells.Item(2, 4) "45.98"
Set-TextValue $q3.Cells.Item(2, 5) "94.53"
Set-TextValue $q3.Cells.Item(2, 6) "2.30"
Set-TextValue $q3.Cells.Item(2, 7) "1.0575"
$q3.Cells.Item(2, 8).Value = 1

Set-TextValue $q3.Cells.Item(3, 3) "华宝油气（QDII）人民币A"
Set-TextValue $q3.Cells.Item(3, 4) "28.25"
Set-TextValue $q3.Cells.Item(3, 5) "94.53"
Set-TextValue $q3.Cells.Item(3, 6) "2.30"
Set-TextValue $q3.Cells.Item(3, 7) "0.6498"
$q3.Cells.Item(3, 8).Value = 1

Set-TextValue $q3.Cells.Item(4, 3) "华宝油气（QDII）人民币 C"
Set-TextValue $q3.Cells.Item(4, 4) "17.73"
Set-TextValue $q3.Cells.Item(4, 5) "94.53"
Set-TextValue $q3.Cells.Item(4, 6) "2.30"
Set-TextValue $q3.Cells.Item(4, 7) "0.4078"
$q3.Cells.Item(4, 8).Value = 1

# ---------------------------------------------------------------------------
# Step 3: update the "总计" (totals) sheet - insert the 2022-Q3 row where the
# 2021-Q4 row used to be, and push the 2021-Q4 row down to row 3 (unchanged
# values).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Copy row 2's formatting down into row 3 first (so the new A3 cell matches
# A2's bold/border style), then fill in the preserved 2021-Q4 figures.
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item(3, 1).PasteSpecial(-4122)  # xlPasteFormats

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 3
$total.Cells.Item(3, 4).Value = 1.63

# Row 2 becomes the 2022-Q3 total.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 2.12

# Leave the view the way it started, focused on the totals sheet.
$total.Activate()

